$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update Riders column (C)
$ws.Range("C3").Value = 178
$ws.Range("C4").Value = 233
$ws.Range("C5").Value = 216
$ws.Range("C6").Value = 120
$ws.Range("C7").Value = 68

# Update Average column (D)
$ws.Range("D2").Value = 218.61
$ws.Range("D3").Value = 215.03
$ws.Range("D4").Value = 234.8
$ws.Range("D5").Value = 237.87
$ws.Range("D6").Value = 114.25
$ws.Range("D7").Value = 90.3
